# Insert a new row at the top of the sheet, shifting the existing header
# row (and all data rows) down by one. The new row 1 becomes a numeric
# index row (0..10); the old header row (now row 2) loses its J/K labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every existing row down by one, inserting a blank row 1.
$ws.Rows("1:1").Insert()

# Carry the header's bold/centered style up onto the new row 1 (it shifted
# down to row 2 along with the content when the row was inserted).
$ws.Range("A2:K2").Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)  # xlPasteFormats

# New row 1: numeric index values 0..10, using the same (bold/header) style
# that the rest of the header row already carries.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10

# The old header row (was row 1, now row 2) keeps its text labels, but the
# last two columns (thread_size / material_surface) are cleared.
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
